$wb = $excel.ActiveWorkbook

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5245.3076
$ws.Range("I100").Value = 6346.4
$ws.Range("J100").Value = 1575
$ws.Range("K100").Value = 6346.4
$ws.Range("L100").Value = 1575
$ws.Range("M100").Value = -5805.4
$ws.Range("N100").Value = -2657

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1408.1666
$ws.Range("I129").Value = 997
$ws.Range("J129").Value = 1490.4
$ws.Range("K129").Value = 2991
$ws.Range("L129").Value = 4471.200000000001
$ws.Range("M129").Value = 2009
$ws.Range("N129").Value = -14471.2

# Sheet ARM, row 64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 49995
$ws.Range("J64").Value = 49995
$ws.Range("L64").Value = 49995
$ws.Range("N64").Value = -50491

# Sheet ARM, row 67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 49995
$ws.Range("J67").Value = 49995
$ws.Range("L67").Value = 49995
$ws.Range("N67").Value = -51711

# Sheet BSM, row 57
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 60000
$ws.Range("J57").Value = 60000
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61440

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2428.5715
$ws.Range("I134").Value = 2333.3333
$ws.Range("K134").Value = 6999.999899999999
$ws.Range("M134").Value = -4464.999899999999

# Sheet BSM, row 136
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 60000
$ws.Range("J136").Value = 60000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

# Sheet BSM, row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3065.6572
$ws.Range("I31").Value = 1466
$ws.Range("J31").Value = 3465.5715
$ws.Range("K31").Value = 1466
$ws.Range("L31").Value = 3465.5715
$ws.Range("M31").Value = -1171
$ws.Range("N31").Value = -4055.5715

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3065.6572
$ws.Range("I34").Value = 1466
$ws.Range("J34").Value = 3465.5715
$ws.Range("K34").Value = 1466
$ws.Range("L34").Value = 3465.5715
$ws.Range("M34").Value = -1264
$ws.Range("N34").Value = -3869.5715

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1971.037
$ws.Range("I107").Value = 1812.381
$ws.Range("K107").Value = 1812.381
$ws.Range("M107").Value = 107.6189999999999

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2789
$ws.Range("I132").Value = 2588.8333
$ws.Range("K132").Value = 7766.499899999999
$ws.Range("M132").Value = -5236.499899999999

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2697.125
$ws.Range("I134").Value = 1746
$ws.Range("K134").Value = 5238
$ws.Range("M134").Value = -2703

# Sheet CRP, row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 110565.445
$ws.Range("J141").Value = 110565.445
$ws.Range("L141").Value = 110565.445
$ws.Range("N141").Value = -120925.445

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18283872
$ws.Range("I4").Value = 52031100
$ws.Range("K4").Value = 156093300
$ws.Range("M4").Value = -156093188

# Sheet CUL, row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 32.75
$ws.Range("I7").Value = 55.5
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 166.5
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = -54.5
$ws.Range("N7").Value = -254

# Sheet CUL, row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 46281.22
$ws.Range("I11").Value = 65216.75
$ws.Range("K11").Value = 195650.25
$ws.Range("M11").Value = -195510.25

# Sheet CUL, row 49
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

# Sheet CUL, row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1079.75
$ws.Range("J97").Value = 606.6
$ws.Range("L97").Value = 1819.8
$ws.Range("N97").Value = -2811.8

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1328
$ws.Range("I113").Value = 683.25
$ws.Range("J113").Value = 1562.4546
$ws.Range("K113").Value = 2049.75
$ws.Range("L113").Value = 4687.3638
$ws.Range("M113").Value = 120.25
$ws.Range("N113").Value = -9027.363799999999

# Sheet CUL, row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2914.1428
$ws.Range("I118").Value = 799.6667
$ws.Range("J118").Value = 4500
$ws.Range("K118").Value = 2399.0001
$ws.Range("L118").Value = 13500
$ws.Range("M118").Value = -1156.0001
$ws.Range("N118").Value = -15986

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1266.3334
$ws.Range("I122").Value = 804
$ws.Range("J122").Value = 1497.5
$ws.Range("K122").Value = 7236
$ws.Range("L122").Value = 13477.5
$ws.Range("M122").Value = -4786
$ws.Range("N122").Value = -18377.5

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1166.1666
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10639.053
$ws.Range("I70").Value = 9546.933999999999
$ws.Range("K70").Value = 9546.933999999999
$ws.Range("M70").Value = -9276.933999999999

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10639.053
$ws.Range("I73").Value = 9546.933999999999
$ws.Range("K73").Value = 9546.933999999999
$ws.Range("M73").Value = -8610.933999999999

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1777.48
$ws.Range("I46").Value = 780.375
$ws.Range("J46").Value = 2246.7058
$ws.Range("K46").Value = 780.375
$ws.Range("L46").Value = 2246.7058
$ws.Range("M46").Value = -592.375
$ws.Range("N46").Value = -2622.7058

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 260.38095
$ws.Range("I55").Value = 231.25
$ws.Range("J55").Value = 353.6
$ws.Range("K55").Value = 231.25
$ws.Range("L55").Value = 353.6
$ws.Range("M55").Value = -58.25
$ws.Range("N55").Value = -699.6

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3090.8276
$ws.Range("I122").Value = 2301.5454
$ws.Range("J122").Value = 5571.4287
$ws.Range("K122").Value = 6904.6362
$ws.Range("L122").Value = 16714.2861
$ws.Range("M122").Value = -4454.6362
$ws.Range("N122").Value = -21614.2861

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2354.1
$ws.Range("I136").Value = 2119.4
$ws.Range("J136").Value = 2588.8
$ws.Range("K136").Value = 6358.200000000001
$ws.Range("L136").Value = 7766.400000000001
$ws.Range("M136").Value = -3808.200000000001

# Sheet WVR, row 44
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 29000
$ws.Range("J44").Value = 29000
$ws.Range("L44").Value = 29000
$ws.Range("N44").Value = -30108

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5957881
$ws.Range("I62").Value = 11906762
$ws.Range("K62").Value = 11906762
$ws.Range("M62").Value = -11906138

# Sheet WVR, row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 28037
$ws.Range("J64").Value = 27056.5
$ws.Range("L64").Value = 27056.5
$ws.Range("N64").Value = -27552.5

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5957881
$ws.Range("I65").Value = 11906762
$ws.Range("K65").Value = 59533810
$ws.Range("M65").Value = -59530690

# Sheet WVR, row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 28037
$ws.Range("J67").Value = 27056.5
$ws.Range("L67").Value = 27056.5
$ws.Range("N67").Value = -28772.5

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 579
$ws.Range("I107").Value = 579
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1737
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 183
$ws.Range("N107").ClearContents()

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5774.6313
$ws.Range("I132").Value = 5976.75
$ws.Range("K132").Value = 17930.25
$ws.Range("M132").Value = -15400.25
